$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.45000000000007
$ws.Range("H2").Value = 0.0000000000000003218037752536686
$ws.Range("K2").Value = 41.21655819288165
$ws.Range("L2").Value = "[32.597647106742585, 49.835469279020714]"
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1.402552876377425
$ws.Range("P2").Value = "[1.1635528346628856, 1.6415529180919641]"
$ws.Range("S2").Value = 58.77879068701205
$ws.Range("T2").Value = "[53.398044840002015, 64.15953653402208]"
$ws.Range("W2").Value = 17.43863863863869
$ws.Range("X2").Value = 16.58468468468473
$ws.Range("Y2").Value = 18.29259259259265

# Row 3 updates
$ws.Range("E3").Value = 22.46000000000007
$ws.Range("H3").Value = 0.0000000000000003218037752536686
$ws.Range("K3").Value = 39.32510796195057
$ws.Range("L3").Value = "[30.702943694213765, 47.94727222968738]"
$ws.Range("O3").Value = 1.314500229429963
$ws.Range("P3").Value = "[1.0755001877154244, 1.5535002711445012]"
$ws.Range("S3").Value = 63.58404808320635
$ws.Range("T3").Value = "[58.63779505080312, 68.53030111560957]"
$ws.Range("W3").Value = 17.76116116116122
$ws.Range("X3").Value = 16.90682682682688
$ws.Range("Y3").Value = 18.61549549549555
